$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-10 Saturday" "2026-01-11 Sunday"

Replace-Text "28×72=2016" "32×11=352"
Replace-Text "95×19=1805" "99×26=2574"
Replace-Text "64×19=1216" "63×14=882"
Replace-Text "42×76=3192" "17×22=374"
Replace-Text "38×11=418" "19×86=1634"

Replace-Text "39×54=2106" "91×69=6279"
Replace-Text "77×81=6237" "56×77=4312"
Replace-Text "72×82=5904" "14×17=238"
Replace-Text "17×49=833" "67×51=3417"
Replace-Text "66×61=4026" "49×39=1911"

Replace-Text "34×15=510" "55×76=4180"
Replace-Text "34×54=1836" "72×72=5184"
Replace-Text "43×71=3053" "16×19=304"
Replace-Text "60×95=5700" "93×74=6882"
Replace-Text "75×43=3225" "30×59=1770"

Replace-Text "81×43=3483" "24×85=2040"
Replace-Text "55×64=3520" "44×92=4048"
Replace-Text "79×88=6952" "16×82=1312"
Replace-Text "55×71=3905" "63×74=4662"
Replace-Text "54×84=4536" "30×19=570"

Replace-Text "26×88=2288" "81×37=2997"
Replace-Text "65×22=1430" "63×74=4662"
Replace-Text "21×54=1134" "65×41=2665"
Replace-Text "32×24=768" "73×31=2263"
Replace-Text "66×12=792" "60×40=2400"
